$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "53.348.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.152.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "398.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.17%  "
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.04%  "
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.645.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.13%  "
$ws.Range("E16").Value = "  +8.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.146.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "53.249.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.55%  "
$ws.Range("E20").Value = "  +5.05%  "
$ws.Range("E21").Value = "  +4.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "271.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.74%  "
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.110"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.03%  "
$ws.Range("E32").Value = "  +7.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0492"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.94%  "
$ws.Range("E40").Value = "  +11.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.294"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  +1.82%  "
$ws.Range("E44").Value = "  +4.35%  "
$ws.Range("E45").Value = "  +1.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.093.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0514"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +29.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0333"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.61%  "